$wb = $excel.ActiveWorkbook

# Fix typo in the "choices" sheet header: list_name -> "list name"
$choices = $wb.Worksheets.Item("choices")
$choices.Range("B2").Value = "list name"

# Move the active selection on "survey" (it stays put at C3, just loses tab focus)
$survey = $wb.Worksheets.Item("survey")
$survey.Range("C3").Select()

# Make "choices" the active/selected sheet, with its selection now at B2
$choices.Activate()
$choices.Range("B2").Select()
